$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 7 ("insertar login al programa"), pushing
#    the existing activities (rows 7-12) down to rows 8-13.
# ---------------------------------------------------------------------------
$ws.Rows(7).Insert()

# ---------------------------------------------------------------------------
# 2) Renumber column A ("No.") for the rows that shifted down one position,
#    since Excel does not auto-increment plain numeric literals on insert.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value  = 3
$ws.Range("A9").Value  = 4
$ws.Range("A10").Value = 5
$ws.Range("A11").Value = 6
$ws.Range("A12").Value = 7
$ws.Range("A13").Value = 8
$ws.Range("A14").Value = 9
$ws.Range("A15").Value = 10

# ---------------------------------------------------------------------------
# 3) Give the brand new row 7 its own sequence number and formatting.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 2
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

# B7 / C7 get the "section title" look (merged, centered, partial borders)
# that matches the other merged activity-name cells such as B9:C9 below.
$ws.Range("B7").Value = "insertar login al programa"
$ws.Range("B9").Copy()
$ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C9").Copy()
$ws.Range("C7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B7:C7").Merge() | Out-Null

# D7 gets a fresh orange highlight fill.
$ws.Range("E6").Copy()
$ws.Range("D7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D7").Interior.Color = 49407   # BGR value for RGB FFC000 (orange)

# E7:M7 stay plain bordered cells, like E6:M6.
$ws.Range("E6:M6").Copy()
$ws.Range("E7:M7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Drop the now-unused trailing blank rows (old rows 16-22, now at 17-23),
#    leaving just two spare rows (14 and 15) after "presentacion".
# ---------------------------------------------------------------------------
$ws.Range("A16:M23").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 5) Refresh the sheet view: drop the frozen top-left cell and point the
#    active selection at the newly blank row (B14:C14).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B14:C14").Select() | Out-Null

Write-Host "edit complete"
